# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de sheets for the first data row.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 12:39:24"
$wsZhCn.Range("H2").Value = "2016-03-22 12:40:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 12:39:32"
$wsDeDe.Range("H2").Value = "2016-03-22 12:40:24"
